$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.020961333333334
$ws.Range("H2").Value = 6.062884
$ws.Range("I2").Value = 0.3447258214530571
$ws.Range("J2").Value = 0.3447258214530571
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.024839333333334
$ws.Range("N2").Value = 9.074518000000001
$ws.Range("O2").Value = 0.1801507982970389
$ws.Range("P2").Value = 0.1801507982970388
$ws.Range("Q2").Value = 6.113083332212446
$ws.Range("R2").Value = 55.01774998991201
$ws.Range("S2").Value = 0.06210263192837072
$ws.Range("T2").Value = 0.06210263192837071

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.020961333333334
$ws.Range("H3").Value = 6.062884
$ws.Range("I3").Value = 0.3447258214530571
$ws.Range("J3").Value = 0.3447258214530571
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.626140333333333
$ws.Range("N3").Value = 4.878420999999999
$ws.Range("O3").Value = 0.09684827751501936
$ws.Range("P3").Value = 0.09684827751501934
$ws.Range("Q3").Value = 3.286366736240444
$ws.Range("R3").Value = 29.577300626164
$ws.Range("S3").Value = 0.03338610202267869
$ws.Range("T3").Value = 0.03338610202267868

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.020961333333334
$ws.Range("H4").Value = 6.062884
$ws.Range("I4").Value = 0.3447258214530571
$ws.Range("J4").Value = 0.3447258214530571
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.640628666666666
$ws.Range("N4").Value = 13.921886
$ws.Range("O4").Value = 0.2763825997921178
$ws.Range("P4").Value = 0.2763825997921177
$ws.Range("Q4").Value = 9.378531097691557
$ws.Range("R4").Value = 84.406779879224
$ws.Range("S4").Value = 0.09527621874866932
$ws.Range("T4").Value = 0.0952762187486693

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.020961333333334
$ws.Range("H5").Value = 6.062884
$ws.Range("I5").Value = 0.3447258214530571
$ws.Range("J5").Value = 0.3447258214530571
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.498988000000001
$ws.Range("N5").Value = 22.496964
$ws.Range("O5").Value = 0.4466183243958241
$ws.Range("P5").Value = 0.446618324395824
$ws.Range("Q5").Value = 15.15516478713067
$ws.Range("R5").Value = 136.396483084176
$ws.Range("S5").Value = 0.1539608687533384
$ws.Range("T5").Value = 0.1539608687533383

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.841556666666667
$ws.Range("H6").Value = 11.52467
$ws.Range("I6").Value = 0.6552741785469429
$ws.Range("J6").Value = 0.6552741785469429
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.024839333333334
$ws.Range("N6").Value = 9.074518000000001
$ws.Range("O6").Value = 0.1801507982970389
$ws.Range("P6").Value = 0.1801507982970388
$ws.Range("Q6").Value = 11.62009170656223
$ws.Range("R6").Value = 104.58082535906
$ws.Range("S6").Value = 0.1180481663686681
$ws.Range("T6").Value = 0.1180481663686681

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.841556666666667
$ws.Range("H7").Value = 11.52467
$ws.Range("I7").Value = 0.6552741785469429
$ws.Range("J7").Value = 0.6552741785469429
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.626140333333333
$ws.Range("N7").Value = 4.878420999999999
$ws.Range("O7").Value = 0.09684827751501936
$ws.Range("P7").Value = 0.09684827751501934
$ws.Range("Q7").Value = 6.246910238452221
$ws.Range("R7").Value = 56.22219214606999
$ws.Range("S7").Value = 0.06346217549234068
$ws.Range("T7").Value = 0.06346217549234066

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.841556666666667
$ws.Range("H8").Value = 11.52467
$ws.Range("I8").Value = 0.6552741785469429
$ws.Range("J8").Value = 0.6552741785469429
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.640628666666666
$ws.Range("N8").Value = 13.921886
$ws.Range("O8").Value = 0.2763825997921178
$ws.Range("P8").Value = 0.2763825997921177
$ws.Range("Q8").Value = 17.82723799195778
$ws.Range("R8").Value = 160.44514192762
$ws.Range("S8").Value = 0.1811063810434485
$ws.Range("T8").Value = 0.1811063810434484

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.841556666666667
$ws.Range("H9").Value = 11.52467
$ws.Range("I9").Value = 0.6552741785469429
$ws.Range("J9").Value = 0.6552741785469429
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.498988000000001
$ws.Range("N9").Value = 22.496964
$ws.Range("O9").Value = 0.4466183243958241
$ws.Range("P9").Value = 0.446618324395824
$ws.Range("Q9").Value = 28.80778734465334
$ws.Range("R9").Value = 259.2700861018801
$ws.Range("S9").Value = 0.2926574556424857
$ws.Range("T9").Value = 0.2926574556424856
